$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell to a text value without letting Excel auto-convert
# numeric-looking strings (e.g. "606.70") into real numbers, and without
# leaving a permanent style change behind on the cell.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.727.30"
$ws.Range("E2").Value = "  +0.64%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.311.47"
$ws.Range("E3").Value = "  +2.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "606.70"
$ws.Range("E5").Value = "  +1.90%  "

# Row 6 - Solana
Set-TextValue "D6" "141.80"
$ws.Range("E6").Value = "  +0.44%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.310.56"
$ws.Range("E8").Value = "  +2.44%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.09%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.67%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +2.92%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.99%  "

# Row 14 - Avalanche
Set-TextValue "D14" "35.01"
$ws.Range("E14").Value = "  +1.94%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.857.37"
$ws.Range("E15").Value = "  +2.43%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +0.37%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.311.96"
$ws.Range("E17").Value = "  +2.48%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "63.798.21"
$ws.Range("E18").Value = "  +0.73%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.86"
$ws.Range("E19").Value = "  +1.13%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "482.02"
$ws.Range("E20").Value = "  +1.82%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -0.84%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +0.67%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "13.98"
$ws.Range("E24").Value = "  +5.98%  "

# Row 25 - Litecoin
Set-TextValue "D25" "85.26"
$ws.Range("E25").Value = "  +1.81%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.01%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  +1.33%  "

# Row 28 - FirstDigitalUSD
$ws.Range("E28").Value = "  -0.02%  "

# Row 29 - RenderToken
$ws.Range("E29").Value = "  +1.46%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  -4.90%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  +1.30%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "28.89"
$ws.Range("E32").Value = "  +5.40%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -0.72%  "

# Row 34 - Stacks
$ws.Range("E34").Value = "  -0.48%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  +1.52%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  +2.47%  "

# Row 37 - OKB
Set-TextValue "D37" "52.49"
$ws.Range("E37").Value = "  -0.25%  "

# Row 38 - PEPE
$ws.Range("E38").Value = "  +5.56%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.0401"
$ws.Range("E39").Value = "  +1.95%  "

# Row 40 - Bittensor
Set-TextValue "D40" "435.02"
$ws.Range("E40").Value = "  +3.01%  "

# Row 41 - Maker
Set-TextValue "D41" "3.113.42"
$ws.Range("E41").Value = "  +4.94%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  +8.59%  "

# Row 43 - Cosmos
Set-TextValue "D43" "8.34"
$ws.Range("E43").Value = "  -0.56%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  -0.39%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -0.26%  "

# Row 46 - Fetch.AI
$ws.Range("E46").Value = "  +3.22%  "

# Row 47 - Arweave
Set-TextValue "D47" "36.85"
$ws.Range("E47").Value = "  +9.48%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "26.36"
$ws.Range("E48").Value = "  +1.51%  "

# Row 50 - ThetaToken
Set-TextValue "D50" "2.32"
$ws.Range("E50").Value = "  -0.99%  "

# Row 51 - Monero
Set-TextValue "D51" "124.63"
$ws.Range("E51").Value = "  +2.98%  "
